# Clustering MALI and Fe# - rescaling + incorporating the time
# Update existing B2 value, and extend the table with rows 3-6,
# reusing the header/label style already applied to column A (style index 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: value rescaled
$ws.Range("B2").Value = 129

# New rows 3-6 (time-incorporated cluster counts)
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 56

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 44

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 39

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 30

# Carry over the formatting used for the existing label column (A2) to the
# newly added label cells A3:A6, so the look & feel stays consistent.
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
